$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "groupedDE_RES" sheet (a dummy RES scenario) right
#    after "powerplants_grouped" and before "germany2019".
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("powerplants_grouped")
$new = $wb.Worksheets.Add()
$new.Name = "groupedDE_RES"
$new.Move($null, $ws3)

# Copy the powerplants_grouped table (header + 11 data rows) as a
# starting point - this carries over the same number formats / fonts /
# borders already used on that sheet.
$ws3.Range("A1:I12").Copy($new.Range("A1"))

# This is a "dummy" scenario: every plant gets a flat age of 8 years
# instead of the real ages from powerplants_grouped.
$new.Range("C2:C12").Value = 8

# The WTG_onshore capacity in the dummy scenario is rounded to 20000.
$new.Range("D12").Value = 20000

# Append extra duplicate RES rows (13-18): two more WTG_onshore @20000
# and four more PV_utility_systems @5000.
$new.Range("A12:I12").Copy($new.Range("A13"))
$new.Range("A12:I12").Copy($new.Range("A14"))
$new.Range("A10:I10").Copy($new.Range("A15"))
$new.Range("A10:I10").Copy($new.Range("A16"))
$new.Range("A10:I10").Copy($new.Range("A17"))
$new.Range("A10:I10").Copy($new.Range("A18"))

$new.Range("A13:A18").Value = 2
$new.Range("D15:D18").Value = 5000

# ------------------------------------------------------------------
# 2. Misc view-state tweaks (selections / zoom / window position) that
#    came along with the edit.
# ------------------------------------------------------------------
$wb.Worksheets.Item("extendedDE").Range("D1:D1048576").Select()

$wsNL = $wb.Worksheets.Item("extendedNL")
$wsNL.Range("D3").Select()
$wsNL.Application.ActiveWindow.Zoom = 130

$wb.Worksheets.Item("powerplants_grouped").Range("B26").Select()

$wb.Worksheets.Item("germany2019").Range("D3").Select()

# Leave the newly added scenario sheet as the active tab/selection.
$new.Range("I27").Select()

$excel.ActiveWindow.Left = 28815
